$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '68.658.75'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.47%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.554.79'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.39%  '

# Row 4
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '594.92'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.72%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '176.76'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.92%  '

# Row 7
$ws.Range('E7').Value = '  -0.06%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.527'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.52%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.554.33'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.41%  '

# Row 10
$ws.Range('E10').Value = '  -0.77%  '

# Row 11
$ws.Range('E11').Value = '  +1.79%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.346'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.03%  '

# Row 13
$ws.Range('E13').Value = '  -2.25%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.66'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.34%  '

# Row 15
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000178'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.55%  '

# Row 16
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.956.98'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.72%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '68.617.74'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.59%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.97'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +98.46%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '2.546.67'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.35%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.98'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +4.76%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '8.06'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.53%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '371.88'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +3.59%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.18'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.13%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.60'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.95%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '72.21'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.29%  '

# Row 26
$ws.Range('E26').Value = '  +0.12%  '

# Row 27
$ws.Range('E27').Value = '  -3.43%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.98'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.63%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0₃0976'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.61%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '539.71'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.83%  '

# Row 32
$ws.Range('E32').Value = '  +0.81%  '

# Row 33
$ws.Range('E33').Value = '  -1.81%  '

# Row 34
$ws.Range('E34').Value = '  +1.33%  '

# Row 35
$ws.Range('E35').Value = '  -0.18%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.12%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '160.39'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.18%  '

# Row 38
$ws.Range('E38').Value = '  -1.50%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '19.35'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +3.31%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.58'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.11%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.19'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.83%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.80'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.83%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.352'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.76%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.55'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.61%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.996'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.44%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '39.47'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.16%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '149.15'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.84%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0₆0281'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.29%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.74'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.16%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.556'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.38%  '

# Row 51
$ws.Range('E51').Value = '  +2.28%  '
